$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 from "N" to "Y" (common utility for runmodes - CustomerSuite now also runs)
$ws.Range("B3").Value = "Y"

# Auto-fit column A width (bestFit) - matches width 17.85546875 in the diff
$ws.Columns.Item(1).AutoFit() | Out-Null

# Move the active cell selection to B2 (as captured in the saved view state)
$ws.Range("B2").Select() | Out-Null
